# Apply cryptos list price/volume updates (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.119.12'
$ws.Range('E2').Value = '  -6.27%  '
$ws.Range('D3').Value = '2.892.28'
$ws.Range('E3').Value = '  -3.53%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '548.05'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -2.62%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '122.33'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -4.38%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '2.887.63'
$ws.Range('E8').Value = '  -3.55%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.497'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('E10').Value = '  -9.65%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '4.64'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -10.74%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.435'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.18%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000210'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -5.81%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '32.43'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.31%  '
$ws.Range('E15').Value = '  +0.90%  '
$ws.Range('D16').Value = '3.360.41'
$ws.Range('E16').Value = '  -3.80%  '
$ws.Range('D17').Value = '2.889.43'
$ws.Range('E17').Value = '  -3.85%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '6.54'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +6.10%  '
$ws.Range('D19').Value = '57.174.50'
$ws.Range('E19').Value = '  -6.30%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '401.92'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -7.45%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.83'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.28%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.670'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.32%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.81'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -4.43%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '12.71'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.15%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '76.99'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -2.32%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.44'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.88%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.92'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.64%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.17'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.10%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '24.64'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -3.29%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.92'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.43%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0991'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +5.71%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.910'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -4.61%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.39'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -3.43%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.99'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -12.19%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '47.80'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -4.54%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '8.24'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +6.05%  '
$ws.Range('D39').Value = '0.0₃0618'
$ws.Range('E39').Value = '  -7.21%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0338'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -5.67%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.105'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -2.28%  '
$ws.Range('D42').Value = '2.625.54'
$ws.Range('E42').Value = '  -2.16%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.39'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.08%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '357.53'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -4.83%  '
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '119.31'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.95%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.228'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -2.71%  '
$ws.Range('E48').Value = '  +0.09%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.93'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.11%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '22.73'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.18%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.94'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -4.01%  '
